$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ETH")
$ws.Range("J3").Value = 3425.049449259545
$ws.Range("B12").Value = 0.00748709
$ws.Range("B36").Value = 0.02575949
$ws.Range("D36").Value = 46.1
$ws.Range("B40").Value = 0.05838402
$ws.Range("D40").Value = 111.45

$ws = $wb.Worksheets.Item("AMP")
$ws.Range("J3").Value = 0.004959247329226439

$ws = $wb.Worksheets.Item("APE")
$ws.Range("J3").Value = 1.951296243019794
$ws.Range("B5").Value = 17.22561887
$ws.Range("D5").Value = 46.1
$ws.Range("B6").Value = 0.60206221

$ws = $wb.Worksheets.Item("ATOM")
$ws.Range("J3").Value = 11.36139464811698
$ws.Range("B7").Value = 0.0299312

$ws = $wb.Worksheets.Item("AVAX")
$ws.Range("J3").Value = 42.20878728470338
$ws.Range("B5").Value = 2.69985508
$ws.Range("D5").Value = 46.1
$ws.Range("B6").Value = 0.0168312

$ws = $wb.Worksheets.Item("BNB")
$ws.Range("J3").Value = 406.5104945119074
$ws.Range("B10").Value = 0.0028247
$ws.Range("B12").Value = 0.15913074
$ws.Range("D12").Value = 46.1

$ws = $wb.Worksheets.Item("DOGE")
$ws.Range("J3").Value = 0.1220978737217196
$ws.Range("B6").Value = 0.29475443

$ws = $wb.Worksheets.Item("DOT")
$ws.Range("J3").Value = 8.458611720656377
$ws.Range("B5").Value = 7.96191652
$ws.Range("D5").Value = 46.1
$ws.Range("B6").Value = 0.08129709

$ws = $wb.Worksheets.Item("EGLD")
$ws.Range("J3").Value = 61.86095533863024
$ws.Range("B6").Value = 0.00300753

$ws = $wb.Worksheets.Item("GRT")
$ws.Range("J3").Value = 0.3138135330356355

$ws = $wb.Worksheets.Item("ICP")
$ws.Range("J3").Value = 12.9556675735905
$ws.Range("B6").Value = 0.00237446

$ws = $wb.Worksheets.Item("BTC")
$ws.Range("J3").Value = 62285.41375382567
$ws.Range("B6").Value = 0.00035713
$ws.Range("B24").Value = 0.00168572
$ws.Range("D24").Value = 46.1
$ws.Range("B34").Value = 0.00219115
$ws.Range("D34").Value = 67.55

$ws = $wb.Worksheets.Item("KAVA")
$ws.Range("J3").Value = 0.8639939127292907

$ws = $wb.Worksheets.Item("LDO")
$ws.Range("J3").Value = 3.401204580373725
$ws.Range("B6").Value = 0.02070408

$ws = $wb.Worksheets.Item("LINK")
$ws.Range("J3").Value = 19.90733330197482
$ws.Range("B6").Value = 0.00250902

$ws = $wb.Worksheets.Item("LTC")
$ws.Range("J3").Value = 84.14605974298716
$ws.Range("B6").Value = 0.00137238

$ws = $wb.Worksheets.Item("LUNA")
$ws.Range("J3").Value = 0.7362680369398577
$ws.Range("B6").Value = 0.0587292

$ws = $wb.Worksheets.Item("LUNC")
$ws.Range("J3").Value = 0.0001430065665303651
$ws.Range("B18").Value = 5075.18440587

$ws = $wb.Worksheets.Item("MATIC")
$ws.Range("J3").Value = 1.018467461558471
$ws.Range("B6").Value = 0.3309549
$ws.Range("B7").Value = 50.55041619
$ws.Range("D7").Value = 46.1

$ws = $wb.Worksheets.Item("MEME")
$ws.Range("J3").Value = 0.03260547835059092
$ws.Range("B6").Value = 0.06972455

$ws = $wb.Worksheets.Item("MINA")
$ws.Range("J3").Value = 1.321764553310866
$ws.Range("B6").Value = 0.35239911

$ws = $wb.Worksheets.Item("NEAR")
$ws.Range("J3").Value = 4.030639191502362
$ws.Range("B6").Value = 24.46551122
$ws.Range("D6").Value = 46.1
$ws.Range("B7").Value = 0.10358297

$ws = $wb.Worksheets.Item("SEI")
$ws.Range("J3").Value = 0.8474442727749686
$ws.Range("B6").Value = 0.07658535

$ws = $wb.Worksheets.Item("SHIB")
$ws.Range("J3").Value = 0.00001318414345048464
$ws.Range("B6").Value = 285.88

$ws = $wb.Worksheets.Item("SHPING")
$ws.Range("J3").Value = 0.006979077712289055

$ws = $wb.Worksheets.Item("SOL")
$ws.Range("J3").Value = 134.4104209373136
$ws.Range("B17").Value = 0.06512804
$ws.Range("B18").Value = 1.93295438
$ws.Range("D18").Value = 46.1

$ws = $wb.Worksheets.Item("TRX")
$ws.Range("J3").Value = 0.1411574570794132
$ws.Range("B6").Value = 0.26978586

$ws = $wb.Worksheets.Item("UNI")
$ws.Range("J3").Value = 11.36187362020049
$ws.Range("B6").Value = 0.00278602

$ws = $wb.Worksheets.Item("XRP")
$ws.Range("J3").Value = 0.5906485519637493
$ws.Range("B6").Value = 0.87892341

$ws = $wb.Worksheets.Item("TIA")
$ws.Range("J3").Value = 16.80012920909715
$ws.Range("B6").Value = 0.00485454

$ws = $wb.Worksheets.Item("DYDX")
$ws.Range("J3").Value = 3.433029405940353
$ws.Range("B6").Value = 0.00102283

$ws = $wb.Worksheets.Item("POLIS")
$ws.Range("J3").Value = 0.4935590871227994

$ws = $wb.Worksheets.Item("ATLAS")
$ws.Range("J3").Value = 0.006847625590008738

$ws = $wb.Worksheets.Item("ACE")
$ws.Range("J3").Value = 12.39979070445509
$ws.Range("B6").Value = 0.00002766

$ws = $wb.Worksheets.Item("ADA")
$ws.Range("J3").Value = 0.685609724819918
$ws.Range("B6").Value = 0.79428926
$ws.Range("B7").Value = 125.66671978
$ws.Range("D7").Value = 46.1

$ws = $wb.Worksheets.Item("ALGO")
$ws.Range("J3").Value = 0.2202802442739022
$ws.Range("B6").Value = 0.58792805
